# Generate Report for Handoff
# - Removes the row for the "b396911e-8a1e-4350-bc5f-2848b741994d" source file
#   (handed-back / in-sync-with-en-US file) from every sheet.
# - Renames the "Handed back: in sync with en-US" status of the remaining
#   "4345290d-..." file to "Ready for handoff".
# - Refreshes its "Latest Handoff Datetime" stamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2 (4345290d file): status text changes in both the zh-cn/de-de columns.
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

# Row 3 is the b396911e-...-md file; drop it entirely. This shifts the
# ".localization-config" row (old row 4) up into row 3, keeping its styles.
$ws.Rows.Item(3).Delete()

# The hyperlink collection does not follow the row shift automatically, so
# rebuild it to match the new two-row layout.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1d5f26a1943622677b6d94b096bcfed8c1677117/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1d5f26a1943622677b6d94b096bcfed8c1677117/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-08 08:30:56"

$ws.Rows.Item(3).Delete()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1d5f26a1943622677b6d94b096bcfed8c1677117/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fcae3ba651b03f512533482b9cdab0500b0ad0ac/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf", "", "", "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f8bd0046d8829fd5b6a4b72adf28f00124f4e8b7/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cb92d7d431986515b31d2ecdf62e429e4830e0b7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf", "", "", "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1d5f26a1943622677b6d94b096bcfed8c1677117/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-08 08:31:00"

$ws.Rows.Item(3).Delete()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1d5f26a1943622677b6d94b096bcfed8c1677117/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f3d1431d29a48a54d57e5ae9138cb82bd6b2f83/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf", "", "", "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e0c3f527d17e0ec73a5f45f61ec84dda8a27d315/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/85c5de4e50bf898c2b0f4e9277c4c0578b4cd11f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf", "", "", "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1d5f26a1943622677b6d94b096bcfed8c1677117/.localization-config", "", "", ".localization-config") | Out-Null

Write-Host "Report regenerated for handoff."
